$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (a new, earlier data point), shifting rows 2-18 down to 3-19
$ws.Rows.Item(2).Insert()

# The new row has no E value (not enough lag yet), so drop that cell entirely
# rather than leaving it as a blank placeholder.
$ws.Range("E2").Clear()

# The insert copies the header row's bold/bordered formatting onto the new
# row; restore plain formatting to match the rest of the data rows, then
# copy the date style (s="2") from the row below onto the new A2 cell.
$ws.Range("A2:D2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Recomputed data for every data row (dates, years, year-over-year growth)
$data = @(
    @(39400, 2007, 0.4235526809466261, 2008, $null),
    @(39765, 2008, -0.571807692896309, 2009, $null),
    @(40130, 2009, 0.3486139762224783, 2010, $null),
    @(40494, 2010, -0.1384957661262898, 2011, 0.592211799485276),
    @(40862, 2011, 1.566479473280147, 2012, 0.5295895589954247),
    @(41228, 2012, 0.7307568962936939, 2013, 0.5955791956549161),
    @(41592, 2013, 0.818818812164257, 2014, 0.9193568360546411),
    @(41957, 2014, 0.9180054319587239, 2015, 0.9697679806505821),
    @(42321, 2015, 1.984684278296656, 2016, 0.970573649360662),
    @(42689, 2016, 1.755995812646982, 2017, 1.083941060573212),
    @(43053, 2017, 1.946965557828384, 2018, 1.343559319682996),
    @(43418, 2018, 1.06432145354225, 2019, 0.8791375467670726),
    @(43783, 2019, 1.361817904277696, 2020, 1.226035857429442),
    @(44159, 2020, -4.352425014431304, 2021, -1.242805832929039),
    @(44525, 2021, -1.761645650979182, 2022, 3.46909983288044),
    @(44890, 2022, 5.20787683103745, 2023, 1.517423464826884),
    @(45254, 2023, -0.9008525709169546, 2024, 0.2090788898015949),
    @(45618, 2024, 0.2738544794132824, 2025, 0.5362040463673612)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    if ($null -ne $entry[2]) {
        $ws.Cells.Item($row, 3).Value = $entry[2]
    }
    $ws.Cells.Item($row, 4).Value = $entry[3]
    if ($null -ne $entry[4]) {
        $ws.Cells.Item($row, 5).Value = $entry[4]
    }
    $row = $row + 1
}
